# Applies the numeric/result updates described by the commit
# "MIP funcionando, tocando CVX en 4node" across the 4node_spain output workbook.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# s_level
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("s_level")
$ws.Range("B2").Value = 1.35427148362301
$ws.Range("C2").Value = 0.5399220128804783
$ws.Range("D2").Value = 0

# ---------------------------------------------------------------------------
# a_level
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("a_level")
$ws.Range("E3").Value = 2.257119139371683
$ws.Range("E4").Value = 0.8998700214674639
$ws.Range("C5").Value = 2.257119139371683
$ws.Range("D5").Value = 0.8998700214674639

# ---------------------------------------------------------------------------
# f_level
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("f_level")
$ws.Range("E2").Value = [double]"8.128085025845606e-10"
$ws.Range("E3").Value = 0.2377944270450742
$ws.Range("E4").Value = 0.2511446015642362
$ws.Range("C5").Value = 0.2353394836361445
$ws.Range("D5").Value = 0.2489821233704967

# ---------------------------------------------------------------------------
# fext_level
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("fext_level")
$ws.Range("E2").Value = 0.9999999991871915
$ws.Range("E3").Value = 0.7622055729549257
$ws.Range("E4").Value = 0.7488553984357638
$ws.Range("C5").Value = 0.7646605163638555
$ws.Range("D5").Value = 0.7510178766295033

# ---------------------------------------------------------------------------
# fij_level
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("fij_level")

# row 2: i2,i2,i4 -> i1,i1,i3
$ws.Range("A2").Value = "i1"
$ws.Range("B2").Value = "i1"
$ws.Range("C2").Value = "i3"
$ws.Range("G2").Value = [double]"4.064042512922803e-10"

# row 3: i3,i3,i4 -> i1,i1,i4
$ws.Range("A3").Value = "i1"
$ws.Range("B3").Value = "i1"
$ws.Range("G3").Value = [double]"4.064042512922803e-10"

# row 4: i4,i4,i1 -> i1,i3,i4
$ws.Range("A4").Value = "i1"
$ws.Range("B4").Value = "i3"
$ws.Range("C4").Value = "i4"
$ws.Range("D4").ClearContents()
$ws.Range("G4").Value = [double]"4.064042512922803e-10"

# row 5: i4,i4,i2 -> i2,i2,i4
$ws.Range("A5").Value = "i2"
$ws.Range("B5").Value = "i2"
$ws.Range("C5").Value = "i4"
$ws.Range("E5").ClearContents()
$ws.Range("G5").Value = 1

# row 6: i4,i4,i3 -> i3,i1,i4
$ws.Range("A6").Value = "i3"
$ws.Range("B6").Value = "i1"
$ws.Range("C6").Value = "i4"
$ws.Range("F6").ClearContents()
$ws.Range("G6").Value = [double]"1.423688050672329e-09"

# row 7: (i3,i3) i4 -> i1
$ws.Range("C7").Value = "i1"
$ws.Range("D7").ClearContents()
$ws.Range("G7").Value = [double]"1.423688050672329e-09"

# row 8: i3,i4,i1 -> i3,i3,i4
$ws.Range("B8").Value = "i3"
$ws.Range("C8").Value = "i4"
$ws.Range("D8").ClearContents()
$ws.Range("G8").Value = 0.9999999983211245

# row 9: i4,i4,i1 -> i4,i3,i1
$ws.Range("B9").Value = "i3"
$ws.Range("D9").Value = [double]"-2.678386409898293e-10"

# row 10: i4,i4,i2 -> i4,i4,i1
$ws.Range("C10").Value = "i1"
$ws.Range("D10").Value = [double]"3.133041357879711e-10"
$ws.Range("E10").ClearContents()

# row 11: i4,i4,i3 -> i4,i4,i2
$ws.Range("C11").Value = "i2"
$ws.Range("E11").Value = 1
$ws.Range("F11").ClearContents()

# row 12: i4,i4,i3 (unchanged labels), F12 stays 1, D12 gets a tiny residual
$ws.Range("D12").Value = [double]"-2.678386409898293e-10"

# ---------------------------------------------------------------------------
# mip_opt_gap
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("mip_opt_gap")
$ws.Range("A1").Value = [double]"1.952074411148966e-08"

# ---------------------------------------------------------------------------
# solver_time
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("solver_time")
$ws.Range("A1").Value = 0.9530000000086147
